$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 220
$ws.Range("I2").Value = 632
$ws.Range("J2").Value = 2691
$ws.Range("K2").Value = 12
$ws.Range("L2").Value = 706
$ws.Range("M2").Value = 39
$ws.Range("N2").Value = 456
$ws.Range("O2").Value = 3
$ws.Range("P2").Value = 11
$ws.Range("Q2").Value = 2
$ws.Range("R2").Value = 22
$ws.Range("S2").Value = 319
$ws.Range("T2").Value = 451
$ws.Range("U2").Value = 41
$ws.Range("V2").Value = 4197
$ws.Range("W2").Value = 1
$ws.Range("X2").Value = 4248
$ws.Range("Z2").Value = 67
$ws.Range("AA2").Value = 29
